# teams-position.xlsx: add stats vs Minsk7x
# - Update the standings table (rows 5-18): games/wins/losses/points + goals string
# - Append two new match-day blocks (29/30 March) with their results

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Copy-CellFormat($srcRow, $srcCol, $dstRow, $dstCol) {
    $ws.Cells.Item($srcRow, $srcCol).Copy() | Out-Null
    $ws.Cells.Item($dstRow, $dstCol).PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------------
# 1. Standings table update (rows 5-18, columns D/E/F/G/H)
#    Row -> [Games, Wins, Losses, "Goals For - Goals Against", Points]
# ---------------------------------------------------------------------------
$standings = @{
    5  = @(13, 11, 2,  "885 - 719",  24)   # ISsoft
    6  = @(13, 10, 3,  "956 - 801",  23)   # Эра-Недвижимости плюс
    7  = @(13, 10, 3,  "897 - 724",  23)   # БГУФК
    8  = @(13, 10, 3,  "1053 - 872", 23)   # Грушвиль
    9  = @(13, 9,  4,  "985 - 806",  22)   # ОПЛАТИ
    10 = @(13, 9,  4,  "930 - 865",  22)   # GOLDEN HILL
    11 = @(13, 8,  5,  "977 - 937",  21)   # Mapogo males
    12 = @(13, 7,  6,  "876 - 743",  20)   # SIRIUS
    13 = @(13, 5,  8,  "784 - 851",  18)   # Стрела
    14 = @(13, 4,  9,  "797 - 904",  17)   # VSS
    15 = @(13, 3,  10, "742 - 817",  16)   # Eagles
    16 = @(13, 3,  10, "713 - 1016", 16)   # NORD
    17 = @(13, 2,  11, "735 - 898",  15)   # ЛФК
    18 = @(13, 0,  13, "572 - 949",  13)   # Минск 7х
}

foreach ($row in $standings.Keys) {
    $vals = $standings[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]   # D: Игры
    $ws.Cells.Item($row, 5).Value = $vals[1]   # E: Побед
    $ws.Cells.Item($row, 6).Value = $vals[2]   # F: Поражений
    $ws.Cells.Item($row, 7).Value = $vals[3]   # G: Мячи
    $ws.Cells.Item($row, 8).Value = $vals[4]   # H: Очки
}

# ---------------------------------------------------------------------------
# 2. Append two new match-day blocks after row 128
#    Row 129: date header 2025-03-29 (serial 45745)
#    Rows 130-132: match results
#    Row 133: date header 2025-03-30 (serial 45746)
#    Rows 134-137: match results
# ---------------------------------------------------------------------------

# Row 129 - date header (copy formats from the existing date-header row 120).
# NOTE: order matters for this engine - Merge() must run BEFORE the format
# copy/paste, otherwise the paste re-derives brand-new (but equivalent)
# border/style objects instead of reusing the existing style indices.
$ws.Range("B129:H129").Merge() | Out-Null
Copy-CellFormat 120 2 129 2
for ($c = 3; $c -le 8; $c++) { Copy-CellFormat 120 3 129 $c }
$ws.Cells.Item(129, 2).Value = 45745

$block1 = @(
    "ЛФК - ISsoft 59:67 (16:30, БНТУ)",
    "ОПЛАТИ - Грушвиль 84:80 (18:00, БНТУ)",
    "Eagles - БГУФК 56:62 (19:30, БНТУ)"
)
$r = 130
foreach ($line in $block1) {
    $ws.Range("B$($r):H$($r)").Merge() | Out-Null
    $ws.Rows.Item($r).RowHeight = 19.95
    Copy-CellFormat 121 2 $r 2
    for ($c = 3; $c -le 8; $c++) { Copy-CellFormat 121 3 $r $c }
    $ws.Cells.Item($r, 2).Value = $line
    $r++
}

# Row 133 - date header
$ws.Range("B133:H133").Merge() | Out-Null
Copy-CellFormat 124 2 133 2
for ($c = 3; $c -le 8; $c++) { Copy-CellFormat 124 3 133 $c }
$ws.Cells.Item(133, 2).Value = 45746

$block2 = @(
    "NORD - VSS 67:70 (11:00, БНТУ)",
    "SIRIUS - Минск 7х 78:48 (12:30, БНТУ)",
    "Эра-Недвижимости плюс - Mapogo males 64:81 (14:00, БНТУ)",
    "GOLDEN HILL - Стрела 79:61 (15:30, БНТУ)"
)
$r = 134
foreach ($line in $block2) {
    $ws.Range("B$($r):H$($r)").Merge() | Out-Null
    $ws.Rows.Item($r).RowHeight = 19.95
    Copy-CellFormat 125 2 $r 2
    for ($c = 3; $c -le 8; $c++) { Copy-CellFormat 125 3 $r $c }
    $ws.Cells.Item($r, 2).Value = $line
    $r++
}
